$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Basic Game rubric")
$ws.Activate()
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "klaar op 09/03/'22"
